$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 15.66598653793335
$ws.Cells.Item(2, 3).Value = 22.79587531089783
$ws.Cells.Item(3, 2).Value = 15.1628360748291
$ws.Cells.Item(3, 3).Value = 22.42717409133911
$ws.Cells.Item(4, 2).Value = 15.81090188026428
$ws.Cells.Item(4, 3).Value = 40.25647974014282
$ws.Cells.Item(5, 2).Value = 16.39934968948364
$ws.Cells.Item(5, 3).Value = 40.51713991165161
$ws.Cells.Item(6, 2).Value = 14.35707259178162
$ws.Cells.Item(6, 3).Value = 23.56313562393188
$ws.Cells.Item(7, 2).Value = 14.512291431427
$ws.Cells.Item(7, 3).Value = 13.89236044883728
$ws.Cells.Item(8, 2).Value = 14.49823951721191
$ws.Cells.Item(8, 3).Value = 13.41277050971985
$ws.Cells.Item(9, 2).Value = 15.27636218070984
$ws.Cells.Item(9, 3).Value = 14.6560001373291
$ws.Cells.Item(10, 2).Value = 16.02639675140381
$ws.Cells.Item(10, 3).Value = 12.96995854377747
$ws.Cells.Item(11, 2).Value = 14.99777603149414
$ws.Cells.Item(11, 3).Value = 13.46646118164062
$ws.Cells.Item(12, 2).Value = 15.39285397529602
$ws.Cells.Item(12, 3).Value = 22.7418270111084
$ws.Cells.Item(13, 2).Value = 14.84666705131531
$ws.Cells.Item(13, 3).Value = 12.44655418395996
$ws.Cells.Item(14, 2).Value = 15.23784160614014
$ws.Cells.Item(14, 3).Value = 21.89424705505371
